# Insert a new weekly price record at row 23 (pushing existing rows 23-48
# down to 24-49), for "Hortaliza, Vega Monumental Concepción - Poroto verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts rows 23..48 down to 24..49
# and carries the date-format style from D23 onward automatically.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record's data.
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44664
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112031
$ws.Range("G23").Value = "Poroto verde"
$ws.Range("H23").Value = "Magnum"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 22000
$ws.Range("L23").Value = 24000
$ws.Range("M23").Value = 23000
$ws.Range("N23").Value = "`$/malla 25 kilos"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 920
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
